$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.635.95"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").Value = "1.883.18"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'249.80"
$ws.Range("E5").Value = "  +1.03%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "'0.4755"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "'0.2939"
$ws.Range("E8").Value = "  +1.28%  "

$ws.Range("D9").Value = "'0.06532"
$ws.Range("E9").Value = "  +0.19%  "

$ws.Range("E10").Value = "  +1.70%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.07738"
$ws.Range("E11").Value = "  +0.04%  "

$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").Value = "'97.05"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").Value = "'0.7399"
$ws.Range("E13").Value = "  -0.69%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.882.42"
$ws.Range("E14").Value = "  +0.27%  "

$ws.Range("D15").Value = "'5.267"
$ws.Range("E15").Value = "  +2.90%  "

$ws.Range("D16").Value = "'274.60"
$ws.Range("E16").Value = "  +0.25%  "

$ws.Range("D17").Value = "30.616.70"
$ws.Range("E17").Value = "  +0.56%  "

$ws.Range("D18").Value = "'13.18"
$ws.Range("E18").Value = "  -3.22%  "

$ws.Range("D19").Value = "'0.000007545"
$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").Value = "2.129.72"
$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("D22").Value = "'5.323"
$ws.Range("E22").Value = "  +1.23%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "'6.241"
$ws.Range("E24").Value = "  +1.21%  "

$ws.Range("D25").Value = "'9.224"
$ws.Range("E25").Value = "  -0.67%  "

$ws.Range("D26").Value = "'164.03"
$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("D27").Value = "'18.88"
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").Value = "'1.917"
$ws.Range("E28").Value = "  -2.02%  "

$ws.Range("D29").Value = "'1.344"
$ws.Range("E29").Value = "  -2.03%  "

$ws.Range("D30").Value = "'0.09693"
$ws.Range("E30").Value = "  -3.05%  "

$ws.Range("E31").Value = "  -0.33%  "

$ws.Range("D32").Value = "'4.300"
$ws.Range("E32").Value = "  -0.59%  "

$ws.Range("D33").Value = "'4.152"
$ws.Range("E33").Value = "  +2.21%  "

$ws.Range("D34").Value = "'0.04873"
$ws.Range("E34").Value = "  +2.00%  "

$ws.Range("E35").Value = "  +0.34%  "

$ws.Range("D36").Value = "'0.6997"
$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D38").Value = "'0.01901"
$ws.Range("E38").Value = "  +1.94%  "

$ws.Range("D39").Value = "'2.773"
$ws.Range("E39").Value = "  +1.39%  "

$ws.Range("D40").Value = "'6.330"
$ws.Range("E40").Value = "  -0.39%  "

$ws.Range("D41").Value = "'74.80"
$ws.Range("E41").Value = "  +6.61%  "

$ws.Range("D42").Value = "'2.018"
$ws.Range("E42").Value = "  +4.54%  "

$ws.Range("D43").Value = "'0.4244"
$ws.Range("E43").Value = "  +1.82%  "

$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("D45").Value = "'0.8403"
$ws.Range("E45").Value = "  +0.40%  "

$ws.Range("D46").Value = "'102.89"
$ws.Range("E46").Value = "  +0.20%  "

$ws.Range("D47").Value = "'9.404"
$ws.Range("E47").Value = "  +0.48%  "

$ws.Range("D48").Value = "'7.061"
$ws.Range("E48").Value = "  -0.23%  "

$ws.Range("D49").Value = "'35.65"
$ws.Range("E49").Value = "  +0.92%  "

$ws.Range("D50").Value = "'917.03"
$ws.Range("E50").Value = "  -0.99%  "

$ws.Range("D51").Value = "'0.05734"
$ws.Range("E51").Value = "  +2.18%  "
